$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# Mapping of column letter -> new value, applied identically to rows 2 and 3
$values = @{
    "D"  = 0.183
    "E"  = 0.0672
    "I"  = 0.004589019573332522
    "J"  = 0.004269351027202547
    "K"  = 3.01
    "L"  = 0.1280851063829787
    "U"  = 8.66
    "V"  = 0.6713178294573643
    "W"  = 0.1127340823970037
    "X"  = 0.03214747629219819
    "Y"  = 0.08058660610480554
    "Z"  = 0.7587794773121016
    "AA" = 0.003239495940882633
    "AB" = 0.02744178820473649
    "AC" = -0.02420229226385385
    "AD" = 20.8
    "AE" = 0.1307902001334287
    "AF" = 20.93079020013343
    "AG" = 12.27079020013343
    "AH" = 0.6186905501264608
    "AI" = 0.435778593542175
    "AJ" = 0.4875011909665188
    "AK" = 0.311672438824757
    "AN" = 155.2238805970149
    "AP" = 91.57306119502559
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
